$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 23: 29-Aug-2023 (serial 45167), "a comprobar", 500
$ws.Range("B23").Value = 45167
$ws.Range("C23").Value = "a comprobar"
$ws.Range("D23").Value = 500

# Row 24: 25-Aug-2023 (serial 45163), "4 botellones", -212
$ws.Range("B24").Value = 45163
$ws.Range("C24").Value = "4 botellones"
$ws.Range("D24").Value = -212

# Update the active selection to match the saved view state (J15)
$ws.Range("J15").Select()

$wb.Save()
